# Regenerate merged AHB files
#
# 1. Rename the "_old" / "_new" header suffixes to the concrete AHB
#    version labels "_FV2310" / "_FV2404" (columns A-J and L-U; column K
#    stays "diff").
# 2. Turn the A1:U73 range into a native Excel Table ("Table1") with an
#    AutoFilter, picking up the (already renamed) header row as its
#    column names.
# 3. Freeze the header row (split below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310", "Segmentgruppe_FV2310", "Segment_FV2310", "Datenelement_FV2310", "Segment ID_FV2310",
    "Code_FV2310", "Qualifier_FV2310", "Beschreibung_FV2310", "Bedingungsausdruck_FV2310", "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
    "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Create the table over the full used range; header names are read from
# the (already renamed) first row.
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U73"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# Freeze panes above row 2 (i.e. keep the header row visible).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
